# Insert a new weekly price record for "Terminal Hortofrutícola Agro Chillán - Piña"
# at row 302, shifting the existing rows 302:366 down to 303:367.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 302 and below down by one row, creating a new blank row 302.
$ws.Rows("302:302").Insert()

# Populate the new row 302 with the latest weekly record.
$ws.Range("A302").Value = 7
$ws.Range("B302").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C302").Value = "Ñuble"
$ws.Range("D302").Value = 45173
$ws.Range("E302").Value = 16
$ws.Range("F302").Value = "Fruta"
$ws.Range("G302").Value = 100108
$ws.Range("H302").Value = "Tropicales y subtropicales"
$ws.Range("I302").Value = 100108005
$ws.Range("J302").Value = "Piña"
$ws.Range("K302").Value = "Caramelo"
$ws.Range("L302").Value = "Tercera"
$ws.Range("M302").Value = 50
$ws.Range("N302").Value = 22000
$ws.Range("O302").Value = 22000
$ws.Range("P302").Value = 22000
$ws.Range("Q302").Value = "$/caja 16 unidades"
$ws.Range("R302").Value = "Ecuador"
$ws.Range("S302").Value = 1375
$ws.Range("T302").Value = 16
